$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (model_3_0_23) - shifts dimension to A1:I24
$ws.Rows(25).Delete()

# Update label text for rows 21 and 22 (order of model_3_0_19/20 swapped)
$ws.Range("A21").Value = "model_3_0_19"
$ws.Range("A22").Value = "model_3_0_20"

# Update numeric metric values for rows 2-24 (new metrics)
$ws.Range("B2").Value = -0.1647107098449927
$ws.Range("C2").Value = -0.005779479693214018
$ws.Range("D2").Value = -0.2853512217752687
$ws.Range("E2").Value = -0.03611169873229558
$ws.Range("F2").Value = 1.288991570472717
$ws.Range("G2").Value = 1.115342378616333
$ws.Range("H2").Value = 1.069115042686462
$ws.Range("I2").Value = 1.093588352203369

$ws.Range("B3").Value = -0.07063640319961006
$ws.Range("C3").Value = -0.0860625859927755
$ws.Range("D3").Value = -0.1175514007727243
$ws.Range("E3").Value = -0.01853898406570709
$ws.Range("F3").Value = 1.184879064559937
$ws.Range("G3").Value = 1.204370975494385
$ws.Range("H3").Value = 0.9295443892478943
$ws.Range("I3").Value = 1.075040817260742

$ws.Range("B4").Value = -0.06332998876146156
$ws.Range("C4").Value = -0.0855037133026606
$ws.Range("D4").Value = -0.1159283845638432
$ws.Range("E4").Value = -0.01762626823003255
$ws.Range("F4").Value = 1.176792979240417
$ws.Range("G4").Value = 1.20375120639801
$ws.Range("H4").Value = 0.9281944036483765
$ws.Range("I4").Value = 1.074077486991882

$ws.Range("B5").Value = -0.04566338852383778
$ws.Range("C5").Value = -0.2227161513419678
$ws.Range("D5").Value = 0.01204873131377671
$ws.Range("E5").Value = -0.04648733345325473
$ws.Range("F5").Value = 1.157241225242615
$ws.Range("G5").Value = 1.355910778045654
$ws.Range("H5").Value = 0.8217470645904541
$ws.Range("I5").Value = 1.104539513587952

$ws.Range("B6").Value = -0.04378630359518021
$ws.Range("C6").Value = -0.221557879492102
$ws.Range("D6").Value = 0.01104632729815336
$ws.Range("E6").Value = -0.04621481481900025
$ws.Range("F6").Value = 1.155163884162903
$ws.Range("G6").Value = 1.354626297950745
$ws.Range("H6").Value = 0.8225808143615723
$ws.Range("I6").Value = 1.104251980781555

$ws.Range("B7").Value = -0.02443699019656576
$ws.Range("C7").Value = -0.1945099251368605
$ws.Range("D7").Value = 0.02797917046688536
$ws.Range("E7").Value = -0.02489050753111233
$ws.Range("F7").Value = 1.133749842643738
$ws.Range("G7").Value = 1.324631810188293
$ws.Range("H7").Value = 0.8084965944290161
$ws.Range("I7").Value = 1.081744432449341

$ws.Range("B8").Value = -0.01845161953215002
$ws.Range("C8").Value = -0.1815285243118687
$ws.Range("D8").Value = -0.02880793671690762
$ws.Range("E8").Value = -0.03872931491799725
$ws.Range("F8").Value = 1.127125859260559
$ws.Range("G8").Value = 1.31023645401001
$ws.Range("H8").Value = 0.8557303547859192
$ws.Range("I8").Value = 1.096351146697998

$ws.Range("B9").Value = -0.007778859587981524
$ws.Range("C9").Value = -0.1834864635534794
$ws.Range("D9").Value = -0.02745364209103052
$ws.Range("E9").Value = -0.03931611872393526
$ws.Range("F9").Value = 1.115314245223999
$ws.Range("G9").Value = 1.312407493591309
$ws.Range("H9").Value = 0.8546038866043091
$ws.Range("I9").Value = 1.096970558166504

$ws.Range("B10").Value = -0.006638213803000337
$ws.Range("C10").Value = -0.1827027293282915
$ws.Range("D10").Value = -0.02854626638662228
$ws.Range("E10").Value = -0.03928542247684907
$ws.Range("F10").Value = 1.114051938056946
$ws.Range("G10").Value = 1.311538457870483
$ws.Range("H10").Value = 0.8555126786231995
$ws.Range("I10").Value = 1.096938133239746

$ws.Range("B11").Value = 0.1360182519616768
$ws.Range("C11").Value = -0.04338493420791023
$ws.Range("D11").Value = 0.2606522025418501
$ws.Range("E11").Value = 0.1454556465891161
$ws.Range("F11").Value = 0.9561732411384583
$ws.Range("G11").Value = 1.157044291496277
$ws.Range("H11").Value = 0.6149664521217346
$ws.Range("I11").Value = 0.901948869228363

$ws.Range("B12").Value = 0.1401488132068653
$ws.Range("C12").Value = -0.03709464138193197
$ws.Range("D12").Value = 0.2591734546710451
$ws.Range("E12").Value = 0.1484060795134108
$ws.Range("F12").Value = 0.9516019225120544
$ws.Range("G12").Value = 1.150068759918213
$ws.Range("H12").Value = 0.6161963939666748
$ws.Range("I12").Value = 0.8988347053527832

$ws.Range("B13").Value = 0.2263285090902388
$ws.Range("C13").Value = 0.03434246766202642
$ws.Range("D13").Value = 0.3173412538189461
$ws.Range("E13").Value = 0.2097126596729875
$ws.Range("F13").Value = 0.85622638463974
$ws.Range("G13").Value = 1.070849895477295
$ws.Range("H13").Value = 0.5678142309188843
$ws.Range("I13").Value = 0.8341272473335266

$ws.Range("B14").Value = 0.2725322162451504
$ws.Range("C14").Value = 0.08473485735699915
$ws.Range("D14").Value = 0.3330422136456769
$ws.Range("E14").Value = 0.2435650159255351
$ws.Range("F14").Value = 0.8050925135612488
$ws.Range("G14").Value = 1.014968037605286
$ws.Range("H14").Value = 0.5547546744346619
$ws.Range("I14").Value = 0.7983969449996948

$ws.Range("B15").Value = 0.2978804967513098
$ws.Range("C15").Value = 0.1069453384619754
$ws.Range("D15").Value = 0.3076536447238747
$ws.Range("E15").Value = 0.2465035948481201
$ws.Range("F15").Value = 0.7770394086837769
$ws.Range("G15").Value = 0.9903380274772644
$ws.Range("H15").Value = 0.5758720636367798
$ws.Range("I15").Value = 0.7952953577041626

$ws.Range("B16").Value = 0.3253754829351155
$ws.Range("C16").Value = 0.1439573106635069
$ws.Range("D16").Value = 0.2820039701734876
$ws.Range("E16").Value = 0.2575783250659095
$ws.Range("F16").Value = 0.7466105222702026
$ws.Range("G16").Value = 0.9492942094802856
$ws.Range("H16").Value = 0.597206711769104
$ws.Range("I16").Value = 0.7836062908172607

$ws.Range("B17").Value = 0.3341185723891006
$ws.Range("C17").Value = 0.1676135971113576
$ws.Range("D17").Value = 0.09020687606382449
$ws.Range("E17").Value = 0.199609403548295
$ws.Range("F17").Value = 0.7369345426559448
$ws.Range("G17").Value = 0.9230610728263855
$ws.Range("H17").Value = 0.7567375898361206
$ws.Range("I17").Value = 0.8447909355163574

$ws.Range("B18").Value = 0.3569374718909658
$ws.Range("C18").Value = 0.1825716295570808
$ws.Range("D18").Value = 0.1415277870593162
$ws.Range("E18").Value = 0.2269617808484019
$ws.Range("F18").Value = 0.7116807699203491
$ws.Range("G18").Value = 0.9064735770225525
$ws.Range("H18").Value = 0.7140504121780396
$ws.Range("I18").Value = 0.8159212470054626

$ws.Range("B19").Value = 0.3645630350584926
$ws.Range("C19").Value = 0.1563492698951623
$ws.Range("D19").Value = 0.1920310904876137
$ws.Range("E19").Value = 0.2311048165272436
$ws.Range("F19").Value = 0.7032414674758911
$ws.Range("G19").Value = 0.9355523586273193
$ws.Range("H19").Value = 0.6720433235168457
$ws.Range("I19").Value = 0.8115484714508057

$ws.Range("B20").Value = 0.3681372046413034
$ws.Range("C20").Value = 0.1583780983873831
$ws.Range("D20").Value = 0.2159048652413301
$ws.Range("E20").Value = 0.2410874280894735
$ws.Range("F20").Value = 0.6992859244346619
$ws.Range("G20").Value = 0.9333025813102722
$ws.Range("H20").Value = 0.65218585729599
$ws.Range("I20").Value = 0.8010119795799255

$ws.Range("B21").Value = 0.3920309924545077
$ws.Range("C21").Value = 0.1911303213986877
$ws.Range("D21").Value = 0.2586953467083167
$ws.Range("E21").Value = 0.2751735099785119
$ws.Range("F21").Value = 0.6728425621986389
$ws.Range("G21").Value = 0.8969824910163879
$ws.Range("H21").Value = 0.6165940761566162
$ws.Range("I21").Value = 0.765035092830658

$ws.Range("B22").Value = 0.4079858277789623
$ws.Range("C22").Value = 0.2141638141926014
$ws.Range("D22").Value = 0.3073269564985193
$ws.Range("E22").Value = 0.3060210502681151
$ws.Range("F22").Value = 0.6551852226257324
$ws.Range("G22").Value = 0.8714399337768555
$ws.Range("H22").Value = 0.5761438608169556
$ws.Range("I22").Value = 0.7324763536453247

$ws.Range("B23").Value = 0.4131129036485381
$ws.Range("C23").Value = 0.2226816719292847
$ws.Range("D23").Value = 0.3074146838691625
$ws.Range("E23").Value = 0.3107902761319483
$ws.Range("F23").Value = 0.6495110392570496
$ws.Range("G23").Value = 0.861994206905365
$ws.Range("H23").Value = 0.5760709047317505
$ws.Range("I23").Value = 0.7274425029754639

$ws.Range("B24").Value = 0.419595252568711
$ws.Range("C24").Value = 0.2396271653841645
$ws.Range("D24").Value = 0.3114535131133291
$ws.Range("E24").Value = 0.3217137661918845
$ws.Range("F24").Value = 0.6423369646072388
$ws.Range("G24").Value = 0.8432028293609619
$ws.Range("H24").Value = 0.5727114677429199
$ws.Range("I24").Value = 0.7159131169319153

